$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 40, pushing the existing rows 40-54 down to 42-56.
$ws.Rows("40:41").Insert()

# Build the accented strings via char codes to avoid any script-encoding issues.
$iacute = [char]0x00ED   # í
$oacute = [char]0x00F3   # ó

$laAraucania = 'La Araucan' + $iacute + 'a'
$espArragos  = 'Esp' + [char]0x00E1 + 'rragos'   # á
$sinEspecificar = 'Sin especificar'
$regionMaule = 'Regi' + $oacute + 'n del Maule'
$regionAraucania = 'Regi' + $oacute + 'n de ' + $laAraucania

# --- Row 40 (new entry) ---
$ws.Cells.Item(40, 1).Value2  = 10
$ws.Cells.Item(40, 2).Value2  = 'Vega Modelo de Temuco'
$ws.Cells.Item(40, 3).Value2  = $laAraucania
$ws.Cells.Item(40, 4).Value2  = 44511
$ws.Cells.Item(40, 5).Value2  = 9
$ws.Cells.Item(40, 6).Value2  = 300000000
$ws.Cells.Item(40, 7).Value2  = $espArragos
$ws.Cells.Item(40, 8).Value2  = $sinEspecificar
$ws.Cells.Item(40, 9).Value2  = 'Primera'
$ws.Cells.Item(40, 10).Value2 = 100
$ws.Cells.Item(40, 11).Value2 = 1300
$ws.Cells.Item(40, 12).Value2 = 1300
$ws.Cells.Item(40, 13).Value2 = 1300
$ws.Cells.Item(40, 14).Value2 = '$/kilo'
$ws.Cells.Item(40, 15).Value2 = $regionAraucania
$ws.Cells.Item(40, 16).Value2 = 1300
$ws.Cells.Item(40, 17).Value2 = 1
$ws.Cells.Item(40, 18).Value2 = 'Hortaliza'

# --- Row 41 (new entry) ---
$ws.Cells.Item(41, 1).Value2  = 10
$ws.Cells.Item(41, 2).Value2  = 'Vega Modelo de Temuco'
$ws.Cells.Item(41, 3).Value2  = $laAraucania
$ws.Cells.Item(41, 4).Value2  = 44511
$ws.Cells.Item(41, 5).Value2  = 9
$ws.Cells.Item(41, 6).Value2  = 300000000
$ws.Cells.Item(41, 7).Value2  = $espArragos
$ws.Cells.Item(41, 8).Value2  = $sinEspecificar
$ws.Cells.Item(41, 9).Value2  = 'Primera'
$ws.Cells.Item(41, 10).Value2 = 350
$ws.Cells.Item(41, 11).Value2 = 1300
$ws.Cells.Item(41, 12).Value2 = 1400
$ws.Cells.Item(41, 13).Value2 = 1357
$ws.Cells.Item(41, 14).Value2 = '$/kilo'
$ws.Cells.Item(41, 15).Value2 = $regionMaule
$ws.Cells.Item(41, 16).Value2 = 1357
$ws.Cells.Item(41, 17).Value2 = 1
$ws.Cells.Item(41, 18).Value2 = 'Hortaliza'
